$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "November Week" -> "Date" (H1)
$ws.Cells.Item(1, 8).Value = "Date"

# Data for columns H (Date), I (Sales), J (Type) for rows 2-70
$rows = @(
    @(2, "9th Nov", 100, "Ducklings"),
    @(3, "9th Nov", 170, "Muscovy Eggs"),
    @(4, "10th Nov", 100, "Ducklings"),
    @(5, "10th Nov", 105, "Chicken Eggs"),
    @(6, "10th Nov", 50, "Pekin Eggs and Hatchery"),
    @(7, "10th Nov", 200, "Ducklings"),
    @(8, "11th Nov", 300, "Ducklings"),
    @(9, "11th Nov", 38.5, "Chicken Eggs"),
    @(10, "12th Nov", 50, "Muscovy Eggs"),
    @(11, "12th Nov", 280, "Chicken Eggs"),
    @(12, "12th Nov", 75, "Muscovy Eggs"),
    @(13, "12th Nov", 105, "Pekin Eggs and Hatchery"),
    @(14, "13th Nov", 210, "Chicken Eggs"),
    @(15, "14th Nov", 24.5, "Chicken Eggs"),
    @(16, "15th Nov", 400, "Ducklings"),
    @(17, "19th Nov", 105, "Guinea Fowl Eggs and Hatchery"),
    @(18, "19th Nov", 360, "Pekin Eggs and Hatchery"),
    @(19, "20th Nov", 45.5, "Chicken Eggs"),
    @(20, "20th Nov", 310, "Ducklings"),
    @(21, "20th Nov", 335, "Chicken Eggs"),
    @(22, "20th Nov", 157.5, "Chicken Eggs"),
    @(23, "20th Nov", 101.5, "Chicken Eggs"),
    @(24, "20th Nov", 40, "Pekin Eggs and Hatchery"),
    @(25, "24th Nov", 320, "Guinea Fowl Chicks"),
    @(26, "27th Nov", 360, "Pekin Eggs and Hatchery"),
    @(27, "27th Nov", 38.5, "Chicken Eggs"),
    @(28, "28th Nov", 560, "Guinea Fowl"),
    @(29, "2nd Dec", 850, "Guinea Fowl Eggs and Hatchery"),
    @(30, "3rd Dec", 100, "Ducklings"),
    @(31, "3rd Dec", 140, "Pekin Eggs and Hatchery"),
    @(32, "4th Dec", 150, "Ducklings"),
    @(33, "4th Dec", 80, "Guinea Fowl Chicks"),
    @(34, "4th Dec", 800, "Guinea Fowl Chicks"),
    @(35, "4th Dec", 100, "Guinea Fowl Eggs and Hatchery"),
    @(36, "4th Dec", 33.5, "Chicken Eggs"),
    @(37, "5th Dec", 420, "Chicken Eggs"),
    @(38, "6th Dec", 400, "Ducklings"),
    @(39, "6th Dec", 160, "Guinea Fowl Chicks"),
    @(40, "6th Dec", 150, "Ducklings"),
    @(41, "6th Dec", 600, "Pekin Eggs and Hatchery"),
    @(42, "8th Dec", 140, "Pekin Eggs and Hatchery"),
    @(43, "8th Dec", 200, "Ducklings"),
    @(44, "9th Dec", 87.5, "Chicken Eggs"),
    @(45, "10th Dec", 45, "Pekin Eggs and Hatchery"),
    @(46, "11th Dec", 100, "Ducklings"),
    @(47, "12th Dec", 160, "Guinea Fowl Chicks"),
    @(48, "12th Dec", 160, "Guinea Fowl Chicks"),
    @(49, "12th Dec", 150, "Guinea Fowl"),
    @(50, "13th Dec", 320, "Guinea Fowl Chicks"),
    @(51, "13th Dec", 100, "Ducklings"),
    @(52, "14th Dec", 31.5, "Chicken Eggs"),
    @(53, "15th Dec", 50, "Ducklings"),
    @(54, "16th Dec", 560, "Guinea Fowl"),
    @(55, "17th Dec", 1000, "Ducklings"),
    @(56, "17th Dec", 300, "Ducklings"),
    @(57, "18th Dec", 100, "Muscovy Eggs"),
    @(58, "19th Dec", 100, "Muscovy Eggs"),
    @(59, "19th Dec", 100, "Ducklings"),
    @(60, "19th Dec", 120, "Muscovy Eggs"),
    @(61, "19th Dec", 160, "Guinea Fowl Chicks"),
    @(62, "20th Dec", 105, "Chicken Eggs"),
    @(63, "21st Dec", 35, "Chicken Eggs"),
    @(64, "21st Dec", 400, "Ducklings"),
    @(65, "24th Dec", 105, "Chicken Eggs"),
    @(66, "25th Dec", 42, "Chicken Eggs"),
    @(67, "25th Dec", 59.5, "Chicken Eggs"),
    @(68, "25th Dec", 94.5, "Chicken Eggs"),
    @(69, "25th Dec", 450, "Guinea Fowl Eggs and Hatchery"),
    @(70, "25th Dec", 50, "Muscovy Eggs")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
}

